# Updated cryptos list with latest price/volume data.
# Applies per-cell content updates to the "Price" (D) and "Volume(1h)" (E)
# columns, plus a reordering of the EnergySwap / Decentraland rows (45/46).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.983.52"
$ws.Range("E2").Value = "  -0.28%  "
$ws.Range("D3").Value = "1.871.67"
$ws.Range("E3").Value = "  -2.60%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "319.48"
$ws.Range("E5").Value = "  -3.72%  "
$ws.Range("E6").Value = "  +0.03%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5088"
$ws.Range("E7").Value = "  -3.23%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3948"
$ws.Range("E8").Value = "  -2.58%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08206"
$ws.Range("E9").Value = "  -3.90%  "
$ws.Range("E10").Value = "  -1.81%  "
$ws.Range("E11").Value = "  -3.09%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "22.93"
$ws.Range("E12").Value = "  +2.12%  "
$ws.Range("D13").Value = "1.865.00"
$ws.Range("E13").Value = "  -2.88%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.292"
$ws.Range("E14").Value = "  -1.93%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.205"
$ws.Range("E15").Value = "  -2.73%  "
$ws.Range("E16").Value = "  +0.15%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "91.88"
$ws.Range("E17").Value = "  -4.98%  "
$ws.Range("E18").Value = "  -2.68%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06377"
$ws.Range("E19").Value = "  -4.89%  "
$ws.Range("E20").Value = "  -2.04%  "
$ws.Range("D22").Value = "29.967.91"
$ws.Range("E22").Value = "  -0.35%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.831"
$ws.Range("E23").Value = "  -3.78%  "
$ws.Range("E24").Value = "  -1.24%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.174"
$ws.Range("E25").Value = "  -2.31%  "
$ws.Range("D26").Value = "2.086.22"
$ws.Range("E26").Value = "  -2.54%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "160.84"
$ws.Range("E27").Value = "  +0.23%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.96"
$ws.Range("E28").Value = "  -0.86%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.243"
$ws.Range("E29").Value = "  -8.93%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "127.51"
$ws.Range("E30").Value = "  -1.71%  "
$ws.Range("E31").Value = "  -1.14%  "
$ws.Range("E32").Value = "  -2.31%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.940"
$ws.Range("E33").Value = "  -3.06%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.730"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.02441"
$ws.Range("E35").Value = "  -3.36%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.221"
$ws.Range("E36").Value = "  +0.27%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06363"
$ws.Range("E37").Value = "  -3.66%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2142"
$ws.Range("E38").Value = "  -4.11%  "
$ws.Range("E39").Value = "  -5.17%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.547"
$ws.Range("E40").Value = "  -5.98%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6306"
$ws.Range("E41").Value = "  -4.07%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "11.31"
$ws.Range("E42").Value = "  -3.20%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.203"
$ws.Range("E43").Value = "  -3.54%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.0000"
$ws.Range("E44").Value = "  +0.08%  "
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "12.98"
$ws.Range("E45").Value = "  -2.58%  "
$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5912"
$ws.Range("E46").Value = "  -4.86%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.642"
$ws.Range("E47").Value = "  -3.95%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.006"
$ws.Range("E48").Value = "  -3.87%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "122.53"
$ws.Range("E49").Value = "  -2.18%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.204"
$ws.Range("E50").Value = "  -3.62%  "
$ws.Range("E51").Value = "  -3.42%  "
